$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving the cell's existing
# number format (prevents Excel from silently re-typing numeric-looking
# text, like "1" or "61.0000", as a real number).
function Set-TextValue {
    param($range, [string]$value)
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# --- Insert a new data row above the current row 7 (CEFTRIAXONE...) ---
# This pushes every existing item row (and the totals/footer rows below)
# down by one, exactly like the target workbook.
$ws.Rows.Item(7).Insert()

# Copy the (now shifted-down) original row 7 formatting from row 8 into the
# freshly inserted, blank row 7 so every cell gets the right style index.
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Restore the row height that PasteSpecial(Formats) does not carry over.
$ws.Rows.Item(7).RowHeight = $ws.Rows.Item(8).RowHeight

# Restore the merged cells for row 7 (PasteSpecial(Formats) doesn't copy
# the merge state, only the per-cell style), matching rows 8-15 below it.
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# --- Populate the new row 7 with the new item's data ---
$ws.Range("A7").Value = 1
Set-TextValue $ws.Range("C7") "BLOKATENS 5/80MG 28 F.C. TAB"
Set-TextValue $ws.Range("H7") "0:1"
Set-TextValue $ws.Range("L7") "1"
Set-TextValue $ws.Range("N7") "122.00"
Set-TextValue $ws.Range("P7") "61.0000"
Set-TextValue $ws.Range("Q7") "0:1"

# --- Update the grand total (now on row 16; was row 15 pre-insert) ---
$ws.Range("P16").Value = 588

# --- Update the generated-at timestamp in the footer (now row 17) ---
$ws.Range("A17").Value = "Sunday, 13 July, 2025 9:41 AM"
